$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.21'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '21.57'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.312'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05640'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.384'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.387'
$ws.Range("E7").Value = '6KuCoinTokenKCS'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8065'
$ws.Range("E8").Value = '7MXTokenMX'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9482'
$ws.Range("E9").Value = '8FTXTokenFTT'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0005841'
$ws.Range("E10").Value = '9OneONE'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1439'
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07438'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03177'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03078'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09261'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.577'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001629'
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.04732'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006383'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.004969'
$ws.Range("E20").Value = '19HotbitTokenHTB'
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001046'
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0001506'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'UpBots'
$ws.Range("C23").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0003111'
$ws.Range("E23").Value = '22UpBotsUBXTWorstin24h'
$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.767'
$ws.Range("E24").Value = '23LEOLEO'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.099'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.3282'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1271'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03954'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006981'
$ws.Range("E41").Value = '40KickTokenKICKBestin24h'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1031'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003082'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007417'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005953'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000753'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0005520'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.6850'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.03308'
$ws.Range("E49").Value = '48BOLOBOLO'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002108'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.01014'
